# Applies the commit "updated based on P4 feedback":
#  1. Bump the fixed datetimeFigureOut field text from 17/05/2016 to
#     27/05/2016 everywhere it appears (the slide master's Date
#     Placeholder plus the same placeholder on every slide layout).
#  2. Merge the "Sem" + "-SOS" runs on slide 7 into a single run
#     "Sem-SOS" (dropping the spell-check err="1" flag that only the
#     first run had).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# 1a. Slide master's Date Placeholder.
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shp = $master.Shapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "27/05/2016"
    }
}

# 1b. Every slide layout's Date Placeholder.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $shapes = $layout.Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "27/05/2016"
        }
    }
}

# 2. Slide 7: merge "Sem" / "-SOS" runs into one "Sem-SOS" run, keeping
#    the second run's formatting (no err="1").
$shape = $p.Slides.Item(7).Shapes.Item(7)
$firstRun = $shape.TextFrame.TextRange.Characters(1, 3)
$firstRun.Text = ""
$shape.TextFrame.TextRange.Text = "Sem-SOS"
